# Apply the balance-tuning edit described in the commit:
# "Minor adjustments to balance out the advantages of each animal type"
#
# - Widen column G (PreyHungerDepletionRate) so its header is fully visible
# - Bump PreyHungerDepletionRate (col G) from 5 to 10 for every data row (2-82)
# - Bump ScentDistance (col J) from 3500 to 4000 for every data row (2-82)
# - Move the active cell selection to E13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G to match the new layout (stored sheet width ends up as 22
# once Excel applies its standard character-width padding of 5/6).
$ws.Columns.Item(7).ColumnWidth = 21.166666666666668

# Update the PreyHungerDepletionRate (G) and ScentDistance (J) values for
# every data row (rows 2 through 82).
for ($row = 2; $row -le 82; $row++) {
    $ws.Cells.Item($row, 7).Value2 = 10
    $ws.Cells.Item($row, 10).Value2 = 4000
}

# Update the selected cell to match the saved view state.
[void]$ws.Range("E13").Select()
